$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the crypto table
# with the latest scraped figures. Values that look like plain decimal
# numbers (e.g. "74.40") are prefixed with a leading apostrophe, exactly as
# typing '74.40 into Excel would do, so the cell keeps the literal text
# (including trailing zeros) instead of being auto-coerced into a number.

$ws.Range("D2").Value = '44.037.59'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '2.360.95'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("E5").Value = '  +1.36%  '
$ws.Range("D6").Value = "`'239.75"
$ws.Range("E6").Value = '  +1.14%  '
$ws.Range("D7").Value = "`'74.40"
$ws.Range("E7").Value = '  +1.71%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = "`'0.598"
$ws.Range("E9").Value = '  +10.71%  '
$ws.Range("D10").Value = "`'0.101"
$ws.Range("E10").Value = '  +1.42%  '
$ws.Range("D11").Value = "`'57.20"
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").Value = "`'32.36"
$ws.Range("E12").Value = '  +12.35%  '
$ws.Range("D13").Value = "`'7.29"
$ws.Range("E13").Value = '  +9.67%  '
$ws.Range("E14").Value = '  +0.84%  '
$ws.Range("D15").Value = '2.712.44'
$ws.Range("E15").Value = '  +0.49%  '
$ws.Range("D16").Value = "`'16.68"
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").Value = "`'0.903"
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").Value = '2.361.58'
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("D19").Value = '43.919.73'
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("E20").Value = '  +0.93%  '
$ws.Range("E21").Value = '  +5.61%  '
$ws.Range("D22").Value = "`'77.06"
$ws.Range("E22").Value = '  -0.89%  '
$ws.Range("D23").Value = "`'257.06"
$ws.Range("E23").Value = '  +1.33%  '
$ws.Range("D24").Value = "`'1.98"
$ws.Range("E24").Value = '  +25.68%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").Value = "`'3.67"
$ws.Range("E26").Value = '  -2.22%  '
$ws.Range("D27").Value = "`'2.50"
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").Value = "`'10.77"
$ws.Range("E28").Value = '  +2.08%  '
$ws.Range("E29").Value = '  -0.72%  '
$ws.Range("D30").Value = "`'22.77"
$ws.Range("E30").Value = '  +1.81%  '
$ws.Range("D31").Value = "`'175.27"
$ws.Range("E31").Value = '  +1.67%  '
$ws.Range("E32").Value = '  -2.35%  '
$ws.Range("E33").Value = '  +3.82%  '
$ws.Range("D34").Value = "`'0.0766"
$ws.Range("E34").Value = '  +7.19%  '
$ws.Range("D35").Value = "`'5.25"
$ws.Range("E35").Value = '  +1.68%  '
$ws.Range("D36").Value = "`'5.46"
$ws.Range("E36").Value = '  +5.33%  '
$ws.Range("E37").Value = '  -6.32%  '
$ws.Range("E38").Value = '  -2.66%  '
$ws.Range("E39").Value = '  -1.18%  '
$ws.Range("E40").Value = '  +4.55%  '
$ws.Range("E41").Value = '  +14.93%  '
$ws.Range("E42").Value = '  +14.06%  '
$ws.Range("E43").Value = '  +4.19%  '
$ws.Range("D44").Value = "`'19.12"
$ws.Range("E44").Value = '  -2.28%  '
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").Value = "`'4.72"
$ws.Range("E46").Value = '  +6.32%  '
$ws.Range("D47").Value = "`'58.76"
$ws.Range("E47").Value = '  +11.88%  '
$ws.Range("E48").Value = '  +7.85%  '
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("E50").Value = '  +1.10%  '
$ws.Range("D51").Value = "`'100.35"
$ws.Range("E51").Value = '  +2.49%  '
